$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The GSC export is a rolling one-day window: the oldest date (2025-11-05,
# row 2) drops off the front and a new date (2026-02-02) is appended after
# the previous last row (row 90). Every other row's Date/Invalid/Valid
# values simply shift up by one row.
#
# Deleting row 2 (a real row delete, not a value overwrite) shifts rows
# 3..90 up to 2..89 intact - cell types/formats move with their cells, so
# the existing text dates keep their shared-string type instead of being
# re-parsed, and no stray styles are introduced.
$ws.Rows.Item(2).Delete()

# Append the new trailing row (now row 90). Format column A as text first
# so Excel stores the ISO date string verbatim instead of auto-converting
# it to a date serial number (matches how the rest of column A is stored).
$ws.Cells.Item(90, 1).NumberFormat = "@"
$ws.Cells.Item(90, 1).Value = "2026-02-02"
$ws.Cells.Item(90, 2).Value = 0
$ws.Cells.Item(90, 3).Value = 28
